$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.563.35"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.451.02"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'578.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'144.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.53%  "
$ws.Range("D7").Value = "3.452.43"
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").Value = "'7.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "4.040.47"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'28.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.13%  "
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.440.30"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "61.720.56"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.97%  "
$ws.Range("D20").Value = "'14.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("D21").Value = "'9.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("D22").Value = "'389.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.68%  "
$ws.Range("D23").Value = "'0.564"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("D24").Value = "'73.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "3.589.05"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "'0.178"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'7.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.26%  "
$ws.Range("E32").Value = "  -9.52%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'24.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").Value = "3.479.60"
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'166.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'28.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.19%  "
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("D44").Value = "'0.801"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'4.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'42.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "2.590.28"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  +2.38%  "
